$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision of the existing A21 timestamp
$ws.Range("A21").Value = 44334.78557194908

# Append the newly retrieved row of data
$ws.Range("A22").Value = 44335.77858406556
$ws.Range("B22").Value = 74107
$ws.Range("C22").Value = 62345
$ws.Range("D22").Value = 3259
$ws.Range("E22").Value = 2051
$ws.Range("F22").Value = 1454
$ws.Range("G22").Value = 19259
$ws.Range("H22").Value = 1389
$ws.Range("I22").Value = 837
$ws.Range("J22").Value = 207
